$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "2019 Hyundai Sante Fe SE"
$ws.Range("B9").Value = 24500
$ws.Range("C9").Value = 29
$ws.Range("D9").Value = "suv"
$ws.Range("E9").Value = "silver"
$ws.Range("F9").Value = "gas"
$ws.Range("G9").Value = 185
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = "all"
$ws.Range("J9").Value = "hyundai"

$ws.Range("A10").Value = "2019 Hyundai Palisade"
$ws.Range("B10").Value = 31550
$ws.Range("C10").Value = 26
$ws.Range("D10").Value = "suv"
$ws.Range("E10").Value = "black"
$ws.Range("F10").Value = "gas"
$ws.Range("G10").Value = 291
$ws.Range("H10").Value = 8
$ws.Range("I10").Value = "all"
$ws.Range("J10").Value = "hyundai"

$ws.Range("J10").Select()
